# Edit Config.xlsx - "Constants" sheet (ActiveSheet / ActiveWorkbook)
# Matches commit: "Calabrio API Process - Added trax report download xaml
# along with all exception handled code"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing file-path values: API_Test -> CalabrioAPI Process ---

# B3: report_filename path (Generate_Reports.xlsm)
$ws.Range("B3").Value = "C:\Users\botfive\Documents\UiPath\CalabrioAPI Process\Generate_Reports.xlsm"

# B6: datewise_folder_path value (Report_data folder)
$ws.Range("B6").Value = "C:\Users\botfive\Documents\UiPath\CalabrioAPI Process\Report_data\"

# B7: report_filepath value (Report_data folder)
$ws.Range("B7").Value = "C:\Users\botfive\Documents\UiPath\CalabrioAPI Process\Report_data\"

# --- Add two new rows of settings below the existing ones ---

# Row 8: trax_URL
$ws.Range("A8").Value = "trax_URL"
$ws.Range("B8").Value = "https://login-eu.calabriocloud.com/?realm=/bravo#/"

# Row 9: Email_id
$ws.Range("A9").Value = "Email_id"
$ws.Range("B9").Value = "vaijayanti.patil@quantanite.com"

# B9 picks up a distinct fill-aware style (applyFill) in the source workbook -
# reproduce the extra cellXfs entry (cellXfs count 8 -> 9, applyFill="1").
$ws.Range("B9").Interior.Color = 16777215

# --- Update the active cell selection on the sheet (now C10) ---
$ws.Range("C10").Select()
